$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '''69.438.09'
$ws.Cells.Item(2, 4).Style = 'Normal'
$ws.Cells.Item(2, 5).Value = '  -1.35%  '
$ws.Cells.Item(2, 5).Style = 'Normal'

$ws.Cells.Item(3, 4).Value = '''2.513.00'
$ws.Cells.Item(3, 4).Style = 'Normal'
$ws.Cells.Item(3, 5).Value = '  -0.35%  '
$ws.Cells.Item(3, 5).Style = 'Normal'

$ws.Cells.Item(4, 4).Value = '''1.00'
$ws.Cells.Item(4, 4).Style = 'Normal'
$ws.Cells.Item(4, 5).Value = '  +0.01%  '
$ws.Cells.Item(4, 5).Style = 'Normal'

$ws.Cells.Item(5, 4).Value = '''572.75'
$ws.Cells.Item(5, 4).Style = 'Normal'
$ws.Cells.Item(5, 5).Value = '  -0.30%  '
$ws.Cells.Item(5, 5).Style = 'Normal'

$ws.Cells.Item(6, 4).Value = '''166.57'
$ws.Cells.Item(6, 4).Style = 'Normal'
$ws.Cells.Item(6, 5).Value = '  -1.71%  '
$ws.Cells.Item(6, 5).Style = 'Normal'

$ws.Cells.Item(7, 5).Value = '  -0.11%  '
$ws.Cells.Item(7, 5).Style = 'Normal'

$ws.Cells.Item(8, 4).Value = '''0.516'
$ws.Cells.Item(8, 4).Style = 'Normal'
$ws.Cells.Item(8, 5).Value = '  +2.19%  '
$ws.Cells.Item(8, 5).Style = 'Normal'

$ws.Cells.Item(9, 4).Value = '''2.514.67'
$ws.Cells.Item(9, 4).Style = 'Normal'
$ws.Cells.Item(9, 5).Value = '  -0.34%  '
$ws.Cells.Item(9, 5).Style = 'Normal'

$ws.Cells.Item(10, 4).Value = '''0.160'
$ws.Cells.Item(10, 4).Style = 'Normal'
$ws.Cells.Item(10, 5).Value = '  -1.30%  '
$ws.Cells.Item(10, 5).Style = 'Normal'

$ws.Cells.Item(11, 5).Value = '  -0.67%  '
$ws.Cells.Item(11, 5).Style = 'Normal'

$ws.Cells.Item(12, 5).Value = '  +4.18%  '
$ws.Cells.Item(12, 5).Style = 'Normal'

$ws.Cells.Item(13, 5).Value = '  +3.68%  '
$ws.Cells.Item(13, 5).Style = 'Normal'

$ws.Cells.Item(14, 4).Value = '''2.981.55'
$ws.Cells.Item(14, 4).Style = 'Normal'
$ws.Cells.Item(14, 5).Value = '  -0.19%  '
$ws.Cells.Item(14, 5).Style = 'Normal'

$ws.Cells.Item(15, 4).Value = '''69.603.02'
$ws.Cells.Item(15, 4).Style = 'Normal'
$ws.Cells.Item(15, 5).Value = '  -1.04%  '
$ws.Cells.Item(15, 5).Style = 'Normal'

$ws.Cells.Item(16, 5).Value = '  -1.85%  '
$ws.Cells.Item(16, 5).Style = 'Normal'

$ws.Cells.Item(17, 4).Value = '''24.86'
$ws.Cells.Item(17, 4).Style = 'Normal'
$ws.Cells.Item(17, 5).Value = '  +0.20%  '
$ws.Cells.Item(17, 5).Style = 'Normal'

$ws.Cells.Item(18, 4).Value = '''2.497.23'
$ws.Cells.Item(18, 4).Style = 'Normal'
$ws.Cells.Item(18, 5).Value = '  -1.52%  '
$ws.Cells.Item(18, 5).Style = 'Normal'

$ws.Cells.Item(19, 4).Value = '''11.34'
$ws.Cells.Item(19, 4).Style = 'Normal'
$ws.Cells.Item(19, 5).Value = '  -1.45%  '
$ws.Cells.Item(19, 5).Style = 'Normal'

$ws.Cells.Item(20, 4).Value = '''7.70'
$ws.Cells.Item(20, 4).Style = 'Normal'
$ws.Cells.Item(20, 5).Value = '  +2.14%  '
$ws.Cells.Item(20, 5).Style = 'Normal'

$ws.Cells.Item(21, 4).Value = '''349.55'
$ws.Cells.Item(21, 4).Style = 'Normal'
$ws.Cells.Item(21, 5).Value = '  -1.78%  '
$ws.Cells.Item(21, 5).Style = 'Normal'

$ws.Cells.Item(22, 5).Value = '  +0.42%  '
$ws.Cells.Item(22, 5).Style = 'Normal'

$ws.Cells.Item(23, 4).Value = '''1.99'
$ws.Cells.Item(23, 4).Style = 'Normal'
$ws.Cells.Item(23, 5).Value = '  +2.04%  '
$ws.Cells.Item(23, 5).Style = 'Normal'

$ws.Cells.Item(24, 5).Value = '  +0.19%  '
$ws.Cells.Item(24, 5).Style = 'Normal'

$ws.Cells.Item(25, 4).Value = '''70.08'
$ws.Cells.Item(25, 4).Style = 'Normal'
$ws.Cells.Item(25, 5).Value = '  +1.48%  '
$ws.Cells.Item(25, 5).Style = 'Normal'

$ws.Cells.Item(26, 4).Value = '''3.99'
$ws.Cells.Item(26, 4).Style = 'Normal'
$ws.Cells.Item(26, 5).Value = '  -1.25%  '
$ws.Cells.Item(26, 5).Style = 'Normal'

$ws.Cells.Item(27, 4).Value = '''8.91'
$ws.Cells.Item(27, 4).Style = 'Normal'
$ws.Cells.Item(27, 5).Value = '  -2.75%  '
$ws.Cells.Item(27, 5).Style = 'Normal'

$ws.Cells.Item(28, 4).Value = '''2.646.45'
$ws.Cells.Item(28, 4).Style = 'Normal'
$ws.Cells.Item(28, 5).Value = '  -0.44%  '
$ws.Cells.Item(28, 5).Style = 'Normal'

$ws.Cells.Item(29, 4).Value = '''1.00'
$ws.Cells.Item(29, 4).Style = 'Normal'
$ws.Cells.Item(29, 5).Value = '  -0.04%  '
$ws.Cells.Item(29, 5).Style = 'Normal'

$ws.Cells.Item(30, 4).Value = '''0.0₃0895'
$ws.Cells.Item(30, 4).Style = 'Normal'
$ws.Cells.Item(30, 5).Value = '  -1.09%  '
$ws.Cells.Item(30, 5).Style = 'Normal'

$ws.Cells.Item(31, 4).Value = '''7.89'
$ws.Cells.Item(31, 4).Style = 'Normal'
$ws.Cells.Item(31, 5).Value = '  +1.17%  '
$ws.Cells.Item(31, 5).Style = 'Normal'

$ws.Cells.Item(32, 4).Value = '''464.17'
$ws.Cells.Item(32, 4).Style = 'Normal'
$ws.Cells.Item(32, 5).Value = '  -2.53%  '
$ws.Cells.Item(32, 5).Style = 'Normal'

$ws.Cells.Item(33, 4).Value = '''1.24'
$ws.Cells.Item(33, 4).Style = 'Normal'
$ws.Cells.Item(33, 5).Value = '  -0.98%  '
$ws.Cells.Item(33, 5).Style = 'Normal'

$ws.Cells.Item(34, 4).Value = '''1.73'
$ws.Cells.Item(34, 4).Style = 'Normal'
$ws.Cells.Item(34, 5).Value = '  -1.00%  '
$ws.Cells.Item(34, 5).Style = 'Normal'

$ws.Cells.Item(35, 5).Value = '  +0.03%  '
$ws.Cells.Item(35, 5).Style = 'Normal'

$ws.Cells.Item(36, 2).Value = 'Kaspa'
$ws.Cells.Item(36, 2).Style = 'Normal'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(36, 3).Style = 'Normal'
$ws.Cells.Item(36, 4).Value = '''0.116'
$ws.Cells.Item(36, 4).Style = 'Normal'
$ws.Cells.Item(36, 5).Value = '  +1.49%  '
$ws.Cells.Item(36, 5).Style = 'Normal'

$ws.Cells.Item(37, 2).Value = 'Monero'
$ws.Cells.Item(37, 2).Style = 'Normal'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(37, 3).Style = 'Normal'
$ws.Cells.Item(37, 4).Value = '''157.45'
$ws.Cells.Item(37, 4).Style = 'Normal'
$ws.Cells.Item(37, 5).Value = '  +0.87%  '
$ws.Cells.Item(37, 5).Style = 'Normal'

$ws.Cells.Item(38, 5).Value = '  +1.10%  '
$ws.Cells.Item(38, 5).Style = 'Normal'

$ws.Cells.Item(39, 4).Value = '''18.51'
$ws.Cells.Item(39, 4).Style = 'Normal'
$ws.Cells.Item(39, 5).Value = '  +0.15%  '
$ws.Cells.Item(39, 5).Style = 'Normal'

$ws.Cells.Item(40, 5).Value = '  -0.04%  '
$ws.Cells.Item(40, 5).Style = 'Normal'

$ws.Cells.Item(41, 4).Value = '''4.76'
$ws.Cells.Item(41, 4).Style = 'Normal'
$ws.Cells.Item(41, 5).Value = '  +1.70%  '
$ws.Cells.Item(41, 5).Style = 'Normal'

$ws.Cells.Item(42, 4).Value = '''0.318'
$ws.Cells.Item(42, 4).Style = 'Normal'
$ws.Cells.Item(42, 5).Value = '  +0.55%  '
$ws.Cells.Item(42, 5).Style = 'Normal'

$ws.Cells.Item(43, 4).Value = '''1.60'
$ws.Cells.Item(43, 4).Style = 'Normal'
$ws.Cells.Item(43, 5).Value = '  -2.08%  '
$ws.Cells.Item(43, 5).Style = 'Normal'

$ws.Cells.Item(44, 4).Value = '''38.33'
$ws.Cells.Item(44, 4).Style = 'Normal'
$ws.Cells.Item(44, 5).Value = '  +0.08%  '
$ws.Cells.Item(44, 5).Style = 'Normal'

$ws.Cells.Item(45, 4).Value = '''2.28'
$ws.Cells.Item(45, 4).Style = 'Normal'
$ws.Cells.Item(45, 5).Value = '  -4.75%  '
$ws.Cells.Item(45, 5).Style = 'Normal'

$ws.Cells.Item(46, 5).Value = '  -12.53%  '
$ws.Cells.Item(46, 5).Style = 'Normal'

$ws.Cells.Item(47, 4).Value = '''141.69'
$ws.Cells.Item(47, 4).Style = 'Normal'
$ws.Cells.Item(47, 5).Value = '  -1.02%  '
$ws.Cells.Item(47, 5).Style = 'Normal'

$ws.Cells.Item(48, 4).Value = '''0.527'
$ws.Cells.Item(48, 4).Style = 'Normal'
$ws.Cells.Item(48, 5).Value = '  +1.12%  '
$ws.Cells.Item(48, 5).Style = 'Normal'

$ws.Cells.Item(49, 4).Value = '''3.49'
$ws.Cells.Item(49, 4).Style = 'Normal'
$ws.Cells.Item(49, 5).Value = '  -0.35%  '
$ws.Cells.Item(49, 5).Style = 'Normal'

$ws.Cells.Item(50, 4).Value = '''0.0730'
$ws.Cells.Item(50, 4).Style = 'Normal'
$ws.Cells.Item(50, 5).Value = '  -0.13%  '
$ws.Cells.Item(50, 5).Style = 'Normal'

$ws.Cells.Item(51, 4).Value = '''1.56'
$ws.Cells.Item(51, 4).Style = 'Normal'
$ws.Cells.Item(51, 5).Value = '  -2.97%  '
$ws.Cells.Item(51, 5).Style = 'Normal'
